$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.264.51'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.22%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.818.62'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.36%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.22%  '
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4654'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.40%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3770'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07410'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8694'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.35%  '
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.821.13'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.679'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.408'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.79%  '
$ws.Range("B15").Value = 'TRON'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07086'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.24%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '92.09'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008763'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.11%  '
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.273.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.311'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.03%  '
$ws.Range("E23").Value = '  +1.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.050.43'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.938'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.39'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.246'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.301'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08931'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7821'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.181'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.516'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.66%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.913'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.19%  '
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.096'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01967'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05248'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.278'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.64%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.370'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +20.76%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5299'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.889'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.02%  '
$ws.Range("E44").Value = '  +0.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.601'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5045'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.44%  '
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '105.46'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.9996'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06326'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.46%  '
